# Word COM-interop script implementing the resume edits described by the diff.
# wdReplaceAll = 2 ; wdFindContinue = 1

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "MISS: $find"
    }
}

# 1. Address line: add missing space after "Gorai-2,"
Replace-Text "Gorai-2,Borivali" "Gorai-2, Borivali"

# 2. Education table: "BE (ThirdYear)" -> "BE (third year)"
Replace-Text "BE (ThirdYear)" "BE (third year)"

# 3. Technical Skills bullet: fix spacing + "Javascript" -> "JavaScript"
Replace-Text "Programming knowledge in C , Python, JAVA, Javascript, HTML5 , CSS , Bootstrap." "Programming knowledge in C, Python, JAVA, JavaScript, HTML5, CSS, Bootstrap."

# 4. "Can work efficiently in WINDOWS and LINUX ." -> remove stray space before period
Replace-Text "LINUX ." "LINUX."

# 5. "DEVOPS" -> "DevOps"
Replace-Text "Deployed projects using DEVOPS tools" "Deployed projects using DevOps tools"

# 6. Typo "conduted" -> "conducted"
Replace-Text "Machine Learning and AI conduted by Amazon." "Machine Learning and AI conducted by Amazon."

# 7. Typo "assosciation" -> "association"
Replace-Text "ATS Learning Solution in assosciation with Microsoft." "ATS Learning Solution in association with Microsoft."

# 8. Udemy courses bullet: tighten comma spacing + "Hands on...in" -> "Hands-on...In" + "Tensorflow" -> "TensorFlow"
Replace-Text "Your First Blockchain” , “Build Responsive Real" "Your First Blockchain”, “Build Responsive Real"
Replace-Text "Hands on Python & R in Data Science" "Hands-on Python & R In Data Science"
Replace-Text "Deep Learning Using Tensorflow”" "Deep Learning Using TensorFlow”"

# 9. Coursera bullet: tighten comma spacing + add second AWS course
Replace-Text "Deep Learning” , “Structured Machine Learning” , “AI for Everyone”" "Deep Learning”, “Structured Machine Learning”, “AI for Everyone”"
Replace-Text "AWS Fundamentals: Going Cloud-Native” by AWS" "AWS Fundamentals: Going Cloud-Native” & “AWS Fundamentals: Migrating to the Cloud” by AWS"

# 10. Internship bullet: "Two month ... as an software" -> "Two-month ... as a software"
Replace-Text "Two month internship" "Two-month internship"
Replace-Text "IT department as an software developer" "IT department as a software developer"

# 11. Typo "Converstional" -> "Conversational"
Replace-Text "Leveraging Converstional AI for Secure Healthcare Assistance." "Leveraging Conversational AI for Secure Healthcare Assistance."

# 12. Add new bullet project after the "INDIA SINGAPORE HACKATHON 2019)." line:
#     "Medical Analytica (A therapy based chatbot for emotion analysis and visualization)"
# Locate the paragraph containing the target text and insert a new paragraph after it,
# inheriting its numbering / list formatting (numId 3).
$newParaText = "Medical Analytica (A therapy based chatbot for emotion analysis and visualization)"
$segoeSplit = 20   # length of "Medical Analytica (A" -> rest gets the Segoe UI styling

$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*INDIA SINGAPORE HACKATHON 2019)*") {
        $insertRange = $p.Range
        $insertRange.Collapse(0)
        $insertRange.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 1)
        $npr = $newPara.Range
        $npr.Text = $newParaText

        # Re-fetch the freshly-typed paragraph range and split off the tail portion
        # ("  therapy based chatbot for emotion analysis and visualization)") into its
        # own Segoe UI / dark-grey colored run, matching the source document styling.
        $full = $d.Paragraphs($i + 1).Range
        $tail = $d.Range($full.Start + $segoeSplit, $full.End - 1)
        $tail.Font.Name = "Segoe UI"
        $tail.Font.NameBi = "Segoe UI"
        $tail.Font.Color = 3025188   # 0x2E2924 little-endian == hex 24292E
        break
    }
}
